# Generate Report for handback
#
# The handback (target) files have now been generated for the zh-cn and
# de-de languages. This:
#   * flips the Status column text from "Ready for handoff" to
#     "Handed back: in sync with en-US" (also reflected on the Overview
#     sheet, which shows the very same status text),
#   * fills in the "Latest Target File" (E) / "Latest Handback File" (F)
#     columns with the same file references as the source (A) / handoff (C)
#     files, turning them into hyperlinks just like A and C,
#   * updates the "Latest Handback DateTime" (G) with the real timestamp
#     instead of the zero-date placeholder.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# The "Overview" sheet's B2/C2/B3/C3 cells hold the very same status text as
# the per-language sheets' Status column, so update them as well.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$addrA2 = ""
$addrC2 = ""
$addrA3 = ""
$addrC3 = ""
$targetA2 = $ws.Range("A2").Address()
$targetC2 = $ws.Range("C2").Address()
$targetA3 = $ws.Range("A3").Address()
$targetC3 = $ws.Range("C3").Address()
foreach ($h in $ws.Hyperlinks) {
    $hAddr = $h.Range.Address()
    if ($hAddr -eq $targetA2) { $addrA2 = $h.Address }
    if ($hAddr -eq $targetC2) { $addrC2 = $h.Address }
    if ($hAddr -eq $targetA3) { $addrA3 = $h.Address }
    if ($hAddr -eq $targetC3) { $addrC3 = $h.Address }
}

$textA2 = $ws.Range("A2").Value2
$textC2 = $ws.Range("C2").Value2
$textA3 = $ws.Range("A3").Value2
$textC3 = $ws.Range("C3").Value2

$ws.Hyperlinks.Add($ws.Range("E2"), $addrA2, "", "", $textA2)
$ws.Hyperlinks.Add($ws.Range("F2"), $addrC2, "", "", $textC2)
$ws.Hyperlinks.Add($ws.Range("E3"), $addrA3, "", "", $textA3)
$ws.Hyperlinks.Add($ws.Range("F3"), $addrC3, "", "", $textC3)

$ws.Range("G2").Value = "2016-01-26 05:15:32"
$ws.Range("G3").Value = "2016-01-26 05:15:32"

# --- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("B2").Value = $newStatus
$ws.Range("B3").Value = $newStatus

$addrA2 = ""
$addrC2 = ""
$addrA3 = ""
$addrC3 = ""
$targetA2 = $ws.Range("A2").Address()
$targetC2 = $ws.Range("C2").Address()
$targetA3 = $ws.Range("A3").Address()
$targetC3 = $ws.Range("C3").Address()
foreach ($h in $ws.Hyperlinks) {
    $hAddr = $h.Range.Address()
    if ($hAddr -eq $targetA2) { $addrA2 = $h.Address }
    if ($hAddr -eq $targetC2) { $addrC2 = $h.Address }
    if ($hAddr -eq $targetA3) { $addrA3 = $h.Address }
    if ($hAddr -eq $targetC3) { $addrC3 = $h.Address }
}

$textA2 = $ws.Range("A2").Value2
$textC2 = $ws.Range("C2").Value2
$textA3 = $ws.Range("A3").Value2
$textC3 = $ws.Range("C3").Value2

$ws.Hyperlinks.Add($ws.Range("E2"), $addrA2, "", "", $textA2)
$ws.Hyperlinks.Add($ws.Range("F2"), $addrC2, "", "", $textC2)
$ws.Hyperlinks.Add($ws.Range("E3"), $addrA3, "", "", $textA3)
$ws.Hyperlinks.Add($ws.Range("F3"), $addrC3, "", "", $textC3)

$ws.Range("G2").Value = "2016-01-26 05:15:50"
$ws.Range("G3").Value = "2016-01-26 05:15:50"
